# 10. / 11.02. und Teil 1 12.02.2025
#
# Appends the worked ISBN-13 / EAN check-digit calculation that was added
# after the previous "131313131313" example, and normalizes the "Normal"
# paragraph style's defaults (no forced hyphen suppression switch,
# explicit 0/0 paragraph spacing and left/"start" justification) the way
# the rest of the document already uses them.

$d = $word.ActiveDocument

# --- 1. Append the new worked example paragraphs at the end of the body ---
$newParagraphTexts = @(
    "",
    "978366263938?",
    "x1 x3 x1 etc =148",
    "148/10 = 14,8",
    "Prüfziffer: 8"
)

foreach ($paraText in $newParagraphTexts) {
    $lastPara = $d.Paragraphs.Last
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    if ($paraText -ne "") {
        $newPara.Range.Text = $paraText
    }
}

# --- 2. Normalize the "Normal" style's default paragraph formatting ---
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.Hyphenation = $false
$normal.ParagraphFormat.SpaceBefore = 0
$normal.ParagraphFormat.SpaceAfter = 0
$normal.ParagraphFormat.Alignment = 0
